$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 265
$ws.Range("C3").Value = 163277
$ws.Range("C4").Value = 154270
$ws.Range("C5").Value = 9008
$ws.Range("C8").Value = 64.68000000000001
